# Update TPM-derived LR-pair metrics for Rtn4-Lingo1 (rows 2-13, cols G-T)
# to reflect the new TPM values used by the scoring scripts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 80.60050966666667
$ws.Cells.Item(2, 8).Value = 241.801529
$ws.Cells.Item(2, 9).Value = 0.17420496858261
$ws.Cells.Item(2, 10).Value = 0.17420496858261
$ws.Cells.Item(2, 15).Value = 0.9349445792302935
$ws.Cells.Item(2, 16).Value = 0.9349445792302935
$ws.Cells.Item(2, 17).Value = 52.80265662395144
$ws.Cells.Item(2, 18).Value = 475.223909615563
$ws.Cells.Item(2, 19).Value = 0.1628719910512948
$ws.Cells.Item(2, 20).Value = 0.1628719910512948
$ws.Cells.Item(3, 7).Value = 80.60050966666667
$ws.Cells.Item(3, 8).Value = 241.801529
$ws.Cells.Item(3, 9).Value = 0.17420496858261
$ws.Cells.Item(3, 10).Value = 0.17420496858261
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.04558433333333334
$ws.Cells.Item(3, 14).Value = 0.136753
$ws.Cells.Item(3, 15).Value = 0.0650554207697065
$ws.Cells.Item(3, 16).Value = 0.06505542076970648
$ws.Cells.Item(3, 17).Value = 3.674120499481889
$ws.Cells.Item(3, 18).Value = 33.06708449533701
$ws.Cells.Item(3, 19).Value = 0.01133297753131519
$ws.Cells.Item(3, 20).Value = 0.01133297753131519
$ws.Cells.Item(4, 7).Value = 93.34790299999999
$ws.Cells.Item(4, 9).Value = 0.2017563980255169
$ws.Cells.Item(4, 10).Value = 0.2017563980255169
$ws.Cells.Item(4, 15).Value = 0.9349445792302935
$ws.Cells.Item(4, 16).Value = 0.9349445792302935
$ws.Cells.Item(4, 17).Value = 61.15367370578032
$ws.Cells.Item(4, 18).Value = 550.383063352023
$ws.Cells.Item(4, 19).Value = 0.1886310506589865
$ws.Cells.Item(4, 20).Value = 0.1886310506589865
$ws.Cells.Item(5, 7).Value = 93.34790299999999
$ws.Cells.Item(5, 9).Value = 0.2017563980255169
$ws.Cells.Item(5, 10).Value = 0.2017563980255169
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.04558433333333334
$ws.Cells.Item(5, 14).Value = 0.136753
$ws.Cells.Item(5, 15).Value = 0.0650554207697065
$ws.Cells.Item(5, 16).Value = 0.06505542076970648
$ws.Cells.Item(5, 17).Value = 4.255201926319667
$ws.Cells.Item(5, 18).Value = 38.296817336877
$ws.Cells.Item(5, 19).Value = 0.01312534736653038
$ws.Cells.Item(5, 20).Value = 0.01312534736653038
$ws.Cells.Item(6, 7).Value = 82.28866066666666
$ws.Cells.Item(6, 8).Value = 246.865982
$ws.Cells.Item(6, 9).Value = 0.177853633995942
$ws.Cells.Item(6, 10).Value = 0.177853633995942
$ws.Cells.Item(6, 15).Value = 0.9349445792302935
$ws.Cells.Item(6, 16).Value = 0.9349445792302935
$ws.Cells.Item(6, 17).Value = 53.90859079175043
$ws.Cells.Item(6, 18).Value = 485.1773171257539
$ws.Cells.Item(6, 19).Value = 0.1662832910009146
$ws.Cells.Item(6, 20).Value = 0.1662832910009146
$ws.Cells.Item(7, 7).Value = 82.28866066666666
$ws.Cells.Item(7, 8).Value = 246.865982
$ws.Cells.Item(7, 9).Value = 0.177853633995942
$ws.Cells.Item(7, 10).Value = 0.177853633995942
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.04558433333333334
$ws.Cells.Item(7, 14).Value = 0.136753
$ws.Cells.Item(7, 15).Value = 0.0650554207697065
$ws.Cells.Item(7, 16).Value = 0.06505542076970648
$ws.Cells.Item(7, 17).Value = 3.751073737382889
$ws.Cells.Item(7, 18).Value = 33.759663636446
$ws.Cells.Item(7, 19).Value = 0.01157034299502738
$ws.Cells.Item(7, 20).Value = 0.01157034299502738
$ws.Cells.Item(8, 7).Value = 45.42364
$ws.Cells.Item(8, 8).Value = 136.27092
$ws.Cells.Item(8, 9).Value = 0.09817585288024938
$ws.Cells.Item(8, 10).Value = 0.09817585288024938
$ws.Cells.Item(8, 15).Value = 0.9349445792302935
$ws.Cells.Item(8, 16).Value = 0.9349445792302935
$ws.Cells.Item(8, 17).Value = 29.75773820102667
$ws.Cells.Item(8, 18).Value = 267.8196438092399
$ws.Cells.Item(8, 19).Value = 0.09178898146169996
$ws.Cells.Item(8, 20).Value = 0.09178898146169996
$ws.Cells.Item(9, 7).Value = 45.42364
$ws.Cells.Item(9, 8).Value = 136.27092
$ws.Cells.Item(9, 9).Value = 0.09817585288024938
$ws.Cells.Item(9, 10).Value = 0.09817585288024938
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.04558433333333334
$ws.Cells.Item(9, 14).Value = 0.136753
$ws.Cells.Item(9, 15).Value = 0.0650554207697065
$ws.Cells.Item(9, 16).Value = 0.06505542076970648
$ws.Cells.Item(9, 17).Value = 2.070606346973333
$ws.Cells.Item(9, 18).Value = 18.63545712276
$ws.Cells.Item(9, 19).Value = 0.006386871418549425
$ws.Cells.Item(9, 20).Value = 0.006386871418549423
$ws.Cells.Item(10, 7).Value = 63.60851399999999
$ws.Cells.Item(10, 8).Value = 190.825542
$ws.Cells.Item(10, 9).Value = 0.1374795175462663
$ws.Cells.Item(10, 10).Value = 0.1374795175462663
$ws.Cells.Item(10, 15).Value = 0.9349445792302935
$ws.Cells.Item(10, 16).Value = 0.9349445792302935
$ws.Cells.Item(10, 17).Value = 41.670934054786
$ws.Cells.Item(10, 18).Value = 375.038406493074
$ws.Cells.Item(10, 19).Value = 0.1285357296850777
$ws.Cells.Item(10, 20).Value = 0.1285357296850777
$ws.Cells.Item(11, 7).Value = 63.60851399999999
$ws.Cells.Item(11, 8).Value = 190.825542
$ws.Cells.Item(11, 9).Value = 0.1374795175462663
$ws.Cells.Item(11, 10).Value = 0.1374795175462663
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.04558433333333334
$ws.Cells.Item(11, 14).Value = 0.136753
$ws.Cells.Item(11, 15).Value = 0.0650554207697065
$ws.Cells.Item(11, 16).Value = 0.06505542076970648
$ws.Cells.Item(11, 17).Value = 2.899551705014
$ws.Cells.Item(11, 18).Value = 26.095965345126
$ws.Cells.Item(11, 19).Value = 0.0089437878611886
$ws.Cells.Item(11, 20).Value = 0.008943787861188598
$ws.Cells.Item(12, 7).Value = 97.40706899999999
$ws.Cells.Item(12, 8).Value = 292.221207
$ws.Cells.Item(12, 9).Value = 0.2105296289694155
$ws.Cells.Item(12, 10).Value = 0.2105296289694155
$ws.Cells.Item(12, 15).Value = 0.9349445792302935
$ws.Cells.Item(12, 16).Value = 0.9349445792302935
$ws.Cells.Item(12, 17).Value = 63.812896945981
$ws.Cells.Item(12, 18).Value = 574.316072513829
$ws.Cells.Item(12, 19).Value = 0.19683353537232
$ws.Cells.Item(12, 20).Value = 0.19683353537232
$ws.Cells.Item(13, 7).Value = 97.40706899999999
$ws.Cells.Item(13, 8).Value = 292.221207
$ws.Cells.Item(13, 9).Value = 0.2105296289694155
$ws.Cells.Item(13, 10).Value = 0.2105296289694155
$ws.Cells.Item(13, 11).Value = 1
$ws.Cells.Item(13, 12).Value = 0.3333333333333333
$ws.Cells.Item(13, 13).Value = 0.04558433333333334
$ws.Cells.Item(13, 14).Value = 0.136753
$ws.Cells.Item(13, 15).Value = 0.0650554207697065
$ws.Cells.Item(13, 16).Value = 0.06505542076970648
$ws.Cells.Item(13, 17).Value = 4.440236302319
$ws.Cells.Item(13, 18).Value = 39.96212672087101
$ws.Cells.Item(13, 19).Value = 0.01369609359709551
$ws.Cells.Item(13, 20).Value = 0.01369609359709551
